$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.650.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.96%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.516.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.54%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.73%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.22%  "

# Row 7
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.516"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.65%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.514.66"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.60%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.43%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.168"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.70%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.361"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.27%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.10%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.975.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.52%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.514.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.95%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000177"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.37%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.35%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.521.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.35%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.84%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.55%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "351.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.18%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.84%  "

# Row 23
$ws.Range("E23").Value = "  +0.60%  "

# Row 24
$ws.Range("E24").Value = "  -0.07%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.24%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.59%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.00%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.647.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.42%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.07%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0893"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.45%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.30%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "465.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.31%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.46%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.40%  "

# Row 35
$ws.Range("E35").Value = "  +0.04%  "

# Row 36
$ws.Range("E36").Value = "  +0.46%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "156.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.21%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.09%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.75%  "

# Row 40
$ws.Range("E40").Value = "  -0.03%  "

# Row 41
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.25%  "

# Row 42
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.319"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.01%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.94%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.15%  "

# Row 45
$ws.Range("E45").Value = "  -6.96%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.58%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.97%  "

# Row 48
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.84%  "

# Row 49
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.523"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.18%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0737"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.43%  "

# Row 51
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.40%  "
